# Append a new attendance record (row 4) to the "15_01" sheet,
# mirroring a new face-recognition detection logged at 19:04.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("15_01")

$ws.Cells.Item(4, 1).Value = "Gabriel Taranto"
$ws.Cells.Item(4, 2).Value = "Presente"
$ws.Cells.Item(4, 3).Value = "19:04"
$ws.Cells.Item(4, 4).Value = "ICA"
